$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-189 all move from serial date 45178 (2023-09-09)
# to serial date 45179 (2023-09-10).
$ws.Range("C2:C189").Value = 45179
